{"js": "// Word only ever keeps a single \"_GoBack\" bookmark (it marks the location\n// of the most recent edit), so drop any existing one before we place a new\n// one further down in the document.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n} catch (e) {\n  // Not present - nothing to remove.\n}\nawait context.sync();\n\n// Locate the placeholder text and replace it with the shortened variable\n// name, keeping the original run-level formatting (font, size, etc.).\nconst results = context.document.body.search(\"{supportingDocsList}\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(\"{supportingDocs}\", \"Replace\");\n  await context.sync();\n\n  // Re-find the opening part of the new placeholder so we can drop a\n  // collapsed range right before the closing brace - this is where Word\n  // parks the \"_GoBack\" bookmark (the location of the most recent edit).\n  const openResults = context.document.body.search(\"{supportingDocs\", { matchCase: true });\n  openResults.load(\"items\");\n  await context.sync();\n\n  if (openResults.items.length > 0) {\n    const openRange = openResults.items[0];\n    const endRange = openRange.getRange(\"End\");\n    endRange.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the placeholder text and replace it with the shortened variable name,\n# keeping the original run-level formatting (font, size, etc.).\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"{supportingDocsList}\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    $range.Text = \"{supportingDocs}\"\n\n    # Re-find the opening part of the new placeholder so we can drop a\n    # collapsed range right before the closing brace - this is where Word\n    # parks the \"_GoBack\" bookmark (the location of the most recent edit).\n    $range2 = $d.Content\n    $find2 = $range2.Find\n    $find2.Text = \"{supportingDocs\"\n    $find2.MatchCase = $true\n    $found2 = $find2.Execute()\n\n    if ($found2) {\n        $range2.Collapse(0)   # wdCollapseEnd\n        $d.Bookmarks.Add(\"_GoBack\", $range2)\n    }\n}\n"}
